$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 327
$ws.Range("M2").Value = 38
$ws.Range("R2").Value = 9
$ws.Range("C3").Value = 976
$ws.Range("J3").Value = 377
$ws.Range("M3").Value = 110
$ws.Range("C4").Value = 637
$ws.Range("J4").Value = 22
$ws.Range("M4").Value = 124
$ws.Range("R4").Value = 11
$ws.Range("C5").Value = 310
$ws.Range("M5").Value = 54
$ws.Range("R5").Value = 16
$ws.Range("C6").Value = 152
$ws.Range("J6").Value = 1
$ws.Range("M6").Value = 28
$ws.Range("C7").Value = 952
$ws.Range("J7").Value = 374
$ws.Range("M7").Value = 239
$ws.Range("R7").Value = 17
$ws.Range("C8").Value = 933
$ws.Range("D8").Value = 1
$ws.Range("J8").Value = 219
$ws.Range("M8").Value = 125
$ws.Range("C9").Value = 389
$ws.Range("M9").Value = 86
$ws.Range("C10").Value = 493
$ws.Range("J10").Value = 29
$ws.Range("M10").Value = 79
$ws.Range("R10").Value = 2
$ws.Range("C11").Value = 454
$ws.Range("J11").Value = 0
$ws.Range("M11").Value = 103
$ws.Range("R11").Value = 22
$ws.Range("C12").Value = 397
$ws.Range("M12").Value = 69
$ws.Range("R12").Value = 5
$ws.Range("C13").Value = 112
$ws.Range("J13").Value = 1
$ws.Range("M13").Value = 23
$ws.Range("R13").Value = 1
$ws.Range("C14").Value = 155
$ws.Range("M14").Value = 15
$ws.Range("C15").Value = 737
$ws.Range("J15").Value = 90
$ws.Range("M15").Value = 134
$ws.Range("C16").Value = 883
$ws.Range("J16").Value = 393
$ws.Range("M16").Value = 111
$ws.Range("R16").Value = 25
$ws.Range("C17").Value = 591
$ws.Range("J17").Value = 41
$ws.Range("M17").Value = 116
$ws.Range("C18").Value = 746
$ws.Range("J18").Value = 143
$ws.Range("M18").Value = 135
$ws.Range("C19").Value = 626
$ws.Range("J19").Value = 56
$ws.Range("M19").Value = 103
$ws.Range("C20").Value = 508
$ws.Range("J20").Value = 31
$ws.Range("M20").Value = 78
$ws.Range("C21").Value = 1144
$ws.Range("J21").Value = 518
$ws.Range("M21").Value = 98
$ws.Range("R21").Value = 3
$ws.Range("C22").Value = 630
$ws.Range("J22").Value = 319
$ws.Range("M22").Value = 60
$ws.Range("C23").Value = 429
$ws.Range("J23").Value = 124
$ws.Range("M23").Value = 101
$ws.Range("C24").Value = 388
$ws.Range("M24").Value = 82
$ws.Range("R24").Value = 20
